$d = $word.ActiveDocument

# --- Step 1: locate the "Betreuungsgutscheine" paragraph and replace it with
# itself plus the new IF_verfuegung.PrintSeitenumbruch field run, a page-break
# paragraph, and the ENDIF field paragraph. ---
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $p = $d.Paragraphs($i)
  if ($p.Range.Text -eq "Betreuungsgutscheine" + [char]13) {
    $target = $p
    break
  }
}
if ($target -eq $null) {
  throw "Could not find 'Betreuungsgutscheine' paragraph"
}
$xml1 = '<w:p w:rsidR="00E243C5" w:rsidRDefault="002863F5" w:rsidP="00166257"><w:pPr><w:pStyle w:val="Text"/><w:tabs><w:tab w:val="left" w:pos="4253"/></w:tabs><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="008F1C1D"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Betreuungsgutscheine</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:instrText xml:space="preserve"> DOCVARIABLE  IF_verfuegung.PrintSeitenumbruch  \* MERGEFORMAT </w:instrText></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:bookmarkStart w:id="3" w:name="_GoBack"/><w:bookmarkEnd w:id="3"/></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cs="Arial"/><w:spacing w:val="8"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:br w:type="page"/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Text"/><w:tabs><w:tab w:val="left" w:pos="4253"/></w:tabs><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:lastRenderedPageBreak/><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:instrText xml:space="preserve"> DOCVARIABLE  ENDIF_verfuegung.PrintSeitenumbruch  \* MERGEFORMAT </w:instrText></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p>'
$target.Range.InsertXML($xml1) | Out-Null

# --- Step 2: locate the "Rechtsmittelbelehrung" paragraph (inside the table)
# and replace it with a version that (a) drops the now-redundant _GoBack
# bookmarkStart, and (b) merges the trailing " " + bookmarkEnd + "sie greifbar..."
# runs into a single run. ---
$target2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $p = $d.Paragraphs($i)
  if ($p.Range.Text -like "Rechtsmittelbelehrung:*") {
    $target2 = $p
    break
  }
}
if ($target2 -eq $null) {
  throw "Could not find 'Rechtsmittelbelehrung' paragraph"
}
$xml2 = '<w:p w:rsidR="00166257" w:rsidRDefault="00166257" w:rsidP="00DA58B0"><w:pPr><w:pStyle w:val="Text"/><w:keepNext/><w:keepLines/><w:tabs><w:tab w:val="left" w:pos="4253"/></w:tabs><w:rPr><w:b/></w:rPr></w:pPr><w:r w:rsidRPr="00995646"><w:rPr><w:b/></w:rPr><w:t>Rechtsmittelbelehrung:</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:br/></w:r><w:r w:rsidRPr="003A41EE"><w:t xml:space="preserve">Gegen diese Verfügung </w:t></w:r><w:r><w:t xml:space="preserve">(den Gutschein) </w:t></w:r><w:r w:rsidRPr="003A41EE"><w:t>kann innert 30 Tagen Beschwerde erhoben we</w:t></w:r><w:r w:rsidRPr="003A41EE"><w:t>r</w:t></w:r><w:r w:rsidRPr="003A41EE"><w:t>den. Die Beschwerdefrist kann nicht verlängert werden. Die Beschwerde ist im Doppel de</w:t></w:r><w:r><w:t>r</w:t></w:r><w:r w:rsidRPr="003A41EE"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Direktion für Bildung, Soziales und Sport, Generalsekretariat, Predigergasse 5, Postfach 275</w:t></w:r><w:r w:rsidRPr="003A41EE"><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>3000 Bern 7</w:t></w:r><w:r w:rsidRPr="003A41EE"><w:t>, zuzustellen. Sie muss (a) angeben, welche Entscheidung anstelle der angefochtenen Verfügung beantragt wird; (b) aus welchen Gründen diese andere Entsche</w:t></w:r><w:r w:rsidRPr="003A41EE"><w:t>i</w:t></w:r><w:r w:rsidRPr="003A41EE"><w:t xml:space="preserve">dung verlangt wird, (c) die Unterschrift der </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="003A41EE"><w:t>beschwerdeführenden</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="003A41EE"><w:t xml:space="preserve"> Partei oder der sie vertr</w:t></w:r><w:r w:rsidRPr="003A41EE"><w:t>e</w:t></w:r><w:r w:rsidRPr="003A41EE"><w:t>tenden Person enthalten. Der Beschwerdeschrift beizulegen sind die Beweismittel, soweit</w:t></w:r><w:r w:rsidR="006E40B5" w:rsidRPr="003A41EE"><w:t xml:space="preserve"> sie greifbar sind, und die angefochtene Verfügung.</w:t></w:r></w:p>'
$target2.Range.InsertXML($xml2) | Out-Null

Write-Output "done"
